$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.478.81'
$ws.Range('E2').Value = '  +0.49%  '
$ws.Range('D3').Value = '1.572.27'
$ws.Range('E3').Value = '  +0.37%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('E5').Value = '  -0.06%  '
$ws.Range('D6').Value = "'291.20"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('D7').Value = "'0.3703"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.86%  '
$ws.Range('D8').Value = "'49.86"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.31%  '
$ws.Range('D9').Value = "'0.3381"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.60%  '
$ws.Range('D10').Value = "'0.07535"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.71%  '
$ws.Range('D11').Value = "'1.141"
$ws.Range('D11').Style = 'Normal'
$ws.Range('D13').Value = "'21.24"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.96%  '
$ws.Range('D14').Value = "'6.015"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.55%  '
$ws.Range('D15').Value = "'6.950"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.62%  '
$ws.Range('D16').Value = '1.573.20'
$ws.Range('E16').Value = '  +0.42%  '
$ws.Range('D17').Value = "'0.00001118"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.26%  '
$ws.Range('D18').Value = "'90.60"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.80%  '
$ws.Range('E19').Value = '  +0.21%  '
$ws.Range('E20').Value = '  -0.07%  '
$ws.Range('D21').Value = "'6.303"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.62%  '
$ws.Range('E22').Value = '  -1.04%  '
$ws.Range('D23').Value = "'12.24"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.53%  '
$ws.Range('D24').Value = '22.487.08'
$ws.Range('E24').Value = '  +0.58%  '
$ws.Range('D25').Value = "'2.369"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.27%  '
$ws.Range('D26').Value = "'2.603"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.28%  '
$ws.Range('D27').Value = "'20.04"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.35%  '
$ws.Range('D28').Value = "'149.22"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.12%  '
$ws.Range('D29').Value = "'5.052"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.59%  '
$ws.Range('E30').Value = '  -0.65%  '
$ws.Range('D31').Value = '1.747.35'
$ws.Range('E31').Value = '  +0.43%  '
$ws.Range('D32').Value = "'1.071"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +8.08%  '
$ws.Range('D33').Value = "'6.219"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.81%  '
$ws.Range('D34').Value = "'2.008"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.25%  '
$ws.Range('D35').Value = "'9.754"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.23%  '
$ws.Range('D36').Value = "'0.08346"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.34%  '
$ws.Range('D37').Value = "'0.02483"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.76%  '
$ws.Range('D38').Value = "'1.362"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.14%  '
$ws.Range('D39').Value = "'0.2300"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.35%  '
$ws.Range('D40').Value = "'0.06536"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.16%  '
$ws.Range('D41').Value = "'5.434"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.62%  '
$ws.Range('D42').Value = "'11.31"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('D43').Value = "'0.6219"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.51%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = "'14.09"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.75%  '
$ws.Range('B45').Value = 'Frax'
$ws.Range('C45').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D45').Value = "'1.001"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.03%  '
$ws.Range('D46').Value = "'3.805"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.04%  '
$ws.Range('D47').Value = "'0.5850"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.39%  '
$ws.Range('D48').Value = "'129.20"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.70%  '
$ws.Range('D49').Value = "'2.071"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.42%  '
$ws.Range('E51').Value = '  +0.13%  '
